$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.6863962134059146
$ws.Range("E4").Value = 0.6215543264920883

$ws.Range("B5").Value = 43.98421712937564
$ws.Range("E5").Value = 36.53131776856103

$ws.Range("B6").Value = 0.1829074
$ws.Range("C6").Value = 0.8170925999999999
$ws.Range("E6").Value = 0.8452689
$ws.Range("F6").Value = 0.1547311

$ws.Range("B7").Value = 20.7248214
$ws.Range("C7").Value = 30.3837742
$ws.Range("E7").Value = 30.3965046
$ws.Range("F7").Value = 21.5807964

$ws.Range("B12").Value = 16

$ws.Range("F21").Value = 25

Write-Output "Done applying forecast updates"
